# Fruta / hortaliza, semanal
# Update dates (column D) and volumes/prices (columns M, N, O, P, S)
# for rows 2, 4-6, 8-13 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44435
$ws.Range("M2").Value = 130

# Row 4
$ws.Range("D4").Value = 44417

# Row 5
$ws.Range("D5").Value = 44476

# Row 6
$ws.Range("D6").Value = 44418
$ws.Range("M6").Value = 40

# Row 8
$ws.Range("D8").Value = 44357
$ws.Range("M8").Value = 35
$ws.Range("N8").Value = 1000
$ws.Range("O8").Value = 1000
$ws.Range("P8").Value = 1000
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44424
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 1200
$ws.Range("O9").Value = 1200
$ws.Range("P9").Value = 1200
$ws.Range("S9").Value = 1200

# Row 10
$ws.Range("D10").Value = 44473
$ws.Range("M10").Value = 120

# Row 11
$ws.Range("D11").Value = 44343
$ws.Range("M11").Value = 60

# Row 12
$ws.Range("D12").Value = 44431
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 1300
$ws.Range("O12").Value = 1300
$ws.Range("P12").Value = 1300
$ws.Range("S12").Value = 1300

# Row 13
$ws.Range("D13").Value = 44432
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 1300
$ws.Range("O13").Value = 1300
$ws.Range("P13").Value = 1300
$ws.Range("S13").Value = 1300
